# #47 Fix nutrition list and quantities
# Update the "Final Value" column (D) with corrected quantities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3002.05
    3  = 428.3
    4  = 60.56
    5  = 123.43
    6  = 32.08
    7  = 125.12
    8  = 862.29
    9  = 1213.84
    10 = 3509.07
    11 = 11.77
    12 = 406.49
    13 = 0.8
    14 = 0.83
    15 = 14.72
    16 = 1.39
    17 = 2.19
    18 = 70.31
    19 = 563.7
    20 = 2
    21 = 298.14
    22 = 650.89
    23 = 6.21
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}
